$d = $word.ActiveDocument

$replacements = @(
    @("555÷5=111, 0", "460÷2=230, 0"),
    @("404÷4=101, 0", "117÷5=23, 2"),
    @("589÷8=73, 5", "988÷3=329, 1"),
    @("953÷2=476, 1", "918÷4=229, 2"),
    @("733÷9=81, 4", "950÷7=135, 5"),
    @("591÷6=98, 3", "727÷3=242, 1"),
    @("861÷3=287, 0", "162÷3=54, 0"),
    @("980÷7=140, 0", "511÷9=56, 7"),
    @("217÷6=36, 1", "482÷3=160, 2"),
    @("476÷4=119, 0", "761÷2=380, 1"),
    @("781÷3=260, 1", "428÷9=47, 5"),
    @("853÷3=284, 1", "514÷4=128, 2"),
    @("834÷5=166, 4", "384÷5=76, 4"),
    @("508÷8=63, 4", "379÷2=189, 1"),
    @("180÷7=25, 5", "925÷8=115, 5"),
    @("103÷4=25, 3", "670÷8=83, 6"),
    @("770÷3=256, 2", "368÷5=73, 3"),
    @("731÷3=243, 2", "784÷7=112, 0"),
    @("327÷2=163, 1", "109÷5=21, 4"),
    @("573÷5=114, 3", "700÷4=175, 0"),
    @("181÷9=20, 1", "604÷6=100, 4"),
    @("389÷6=64, 5", "543÷8=67, 7"),
    @("554÷3=184, 2", "750÷9=83, 3"),
    @("907÷7=129, 4", "188÷7=26, 6"),
    @("572÷9=63, 5", "505÷6=84, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
